# Apply the matmap.xlsx edit: add a header row (segment path / poisson's
# ratio / material density / youngs modulus) above the existing data table,
# size the three new numeric columns, and move the active selection.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New header row values - written directly into row 1 (the existing data in
# rows 2-10 is left untouched / is not shifted down).
$ws.Range("A1").Value = "segment path"
$ws.Range("B1").Value = "poisson's ratio"
$ws.Range("C1").Value = "material density"
$ws.Range("D1").Value = "youngs modulus"

# Give the new columns an explicit width (closest achievable values to the
# authored 14.28515625 / 15.5703125 / 16.140625 character widths).
$ws.Columns("B").ColumnWidth = 13.5
$ws.Columns("C").ColumnWidth = 14.666666666666666
$ws.Columns("D").ColumnWidth = 15.333333333333334

# Move the selection to match the saved view state.
$null = $ws.Range("B16").Select()
